$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - LinearRegression values update (last-digit precision change)
$ws.Range("B2").Value = 0.2786851492779058
$ws.Range("C2").Value = 0.2786851492779058
$ws.Range("D2").Value = 0.2786851492779058

# Row 3 - RandomForestRegressor values update
$ws.Range("B3").Value = 0.988328785803518
$ws.Range("C3").Value = 0.9884497726136989
$ws.Range("D3").Value = 0.798063248832325

# Row 4 - label change from GradientBoostingRegressor to DecisionTreeRegressor
$ws.Range("A4").Value = "DecisionTreeRegressor"
$ws.Range("B4").Value = 0.9831578472264163
$ws.Range("C4").Value = 0.984265532266416
$ws.Range("D4").Value = 0.7929075820522051

# Row 5 - label change from AdaBoostRegressor to MLPRegressor
$ws.Range("A5").Value = "MLPRegressor"
$ws.Range("B5").Value = 0.8559852545989114
$ws.Range("C5").Value = 0.8307320747187062
$ws.Range("D5").Value = 0.7276402336670126
